# Update "想去人数" (interest count) figures on the 展览 and 全部类型 sheets,
# matching the regenerated data output.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 12923
$ws1.Range("F10").Value = 12877
$ws1.Range("F13").Value = 8695
$ws1.Range("F14").Value = 7694
$ws1.Range("F22").Value = 383

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 12923
$ws4.Range("F11").Value = 12877
$ws4.Range("F14").Value = 8695
$ws4.Range("F15").Value = 7694
$ws4.Range("F24").Value = 383
